# Apply odds updates to the "Jogos da Semana" worksheet.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 10: Leones del Norte vs Gualaceo - odds were empty, now populated.
$ws.Range("G10").Value = 2.42
$ws.Range("H10").Value = 3.05
$ws.Range("I10").Value = 2.85
$ws.Range("L10").Value = 1.44
$ws.Range("M10").Value = 2.42
$ws.Range("N10").Value = 2.25
$ws.Range("O10").Value = 1.5
$ws.Range("P10").Value = 1.52
$ws.Range("Q10").Value = 2.22
$ws.Range("R10").Value = 1.93
$ws.Range("S10").Value = 1.7
$ws.Range("T10").Value = 6.5
$ws.Range("U10").Value = 10.75
$ws.Range("V10").Value = 10
$ws.Range("W10").Value = 25
$ws.Range("X10").Value = 23
$ws.Range("Y10").Value = 40
$ws.Range("Z10").Value = 7.2
$ws.Range("AA10").Value = 6
$ws.Range("AB10").Value = 17
$ws.Range("AC10").Value = 100
$ws.Range("AE10").Value = 7.2
$ws.Range("AF10").Value = 13
$ws.Range("AG10").Value = 11
$ws.Range("AH10").Value = 35
$ws.Range("AI10").Value = 29
$ws.Range("AJ10").Value = 45

# Row 11: Imbabura vs Nueve de Octubre - odds were empty, now populated.
$ws.Range("G11").Value = 1.95
$ws.Range("H11").Value = 3.25
$ws.Range("I11").Value = 3.7
$ws.Range("L11").Value = 1.36
$ws.Range("M11").Value = 2.65
$ws.Range("N11").Value = 2.05
$ws.Range("O11").Value = 1.6
$ws.Range("P11").Value = 1.47
$ws.Range("Q11").Value = 2.35
$ws.Range("R11").Value = 1.88
$ws.Range("S11").Value = 1.72
$ws.Range("T11").Value = 6.3
$ws.Range("U11").Value = 8.5
$ws.Range("V11").Value = 8.75
$ws.Range("W11").Value = 16.5
$ws.Range("X11").Value = 17
$ws.Range("Y11").Value = 32
$ws.Range("Z11").Value = 8
$ws.Range("AA11").Value = 6.3
$ws.Range("AB11").Value = 16.5
$ws.Range("AC11").Value = 90
$ws.Range("AD11").Value = 800
$ws.Range("AE11").Value = 9.25
$ws.Range("AF11").Value = 19
$ws.Range("AG11").Value = 13
$ws.Range("AH11").Value = 55
$ws.Range("AI11").Value = 40
$ws.Range("AJ11").Value = 50

# Row 13: R. Oviedo vs Almeria - odds adjustments.
$ws.Range("G13").Value = 2.5
$ws.Range("I13").Value = 2.7
$ws.Range("L13").Value = 1.33
$ws.Range("M13").Value = 3.25
$ws.Range("N13").Value = 2.1
$ws.Range("O13").Value = 1.7

# Row 14: Louisville City vs Rhode Island - odds adjustments.
$ws.Range("G14").Value = 1.5
$ws.Range("H14").Value = 3.95
$ws.Range("I14").Value = 5.5
$ws.Range("K14").Value = 7.8
$ws.Range("L14").Value = 1.26
$ws.Range("M14").Value = 3.45
$ws.Range("N14").Value = 1.78
$ws.Range("O14").Value = 1.93
$ws.Range("R14").Value = 1.9
$ws.Range("S14").Value = 1.82
$ws.Range("T14").Value = 6.7
$ws.Range("U14").Value = 7
$ws.Range("W14").Value = 10.5
$ws.Range("X14").Value = 12.5
$ws.Range("Z14").Value = 7.8
$ws.Range("AA14").Value = 7.8
$ws.Range("AE14").Value = 15
$ws.Range("AF14").Value = 35
$ws.Range("AG14").Value = 18
$ws.Range("AH14").Value = 110
$ws.Range("AI14").Value = 60
$ws.Range("AJ14").Value = 60
